$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab07")

# Fix mangled accented characters in the Regional Economic Communities note
# (PALOP / MERCOSUR definitions), cell A103
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# Updated data values throughout the table (rows 13-98)
$ws.Range("F13").Value = 76.02
$ws.Range("G13").Value = 78.56
$ws.Range("H13").Value = 73.52

$ws.Range("G38").Value = 82.29

$ws.Range("F62").Value = 82.03
$ws.Range("G62").Value = 84.427500000000094
$ws.Range("H62").Value = 79.86

$ws.Range("F63").Value = 41.306896551724101
$ws.Range("G63").Value = 40.655172413793103
$ws.Range("H63").Value = 41.704597701149403

$ws.Range("F64").Value = 56.522727272727302
$ws.Range("G64").Value = 54.781818181818203
$ws.Range("H64").Value = 57.818181818181799

$ws.Range("F65").Value = 73.14
$ws.Range("G65").Value = 72.665000000000006
$ws.Range("H65").Value = 72.95

$ws.Range("F66").Value = 54.1330708661418
$ws.Range("G66").Value = 54.441732283464603
$ws.Range("H66").Value = 53.722047244094497

$ws.Range("G67").Value = 80.694117647058903

$ws.Range("I70").Value = 72.221598231208304

$ws.Range("F73").Value = 77.540000000000006
$ws.Range("G73").Value = 79.453333333333404
$ws.Range("H73").Value = 75.62

$ws.Range("F76").Value = 72
$ws.Range("G76").Value = 73
$ws.Range("H76").Value = 71.342857142857198

$ws.Range("F77").Value = 54.127272727272697
$ws.Range("G77").Value = 53.154545454545499
$ws.Range("H77").Value = 54.727272727272698

$ws.Range("F78").Value = 8.0346153846153907
$ws.Range("G78").Value = 7.2923076923077002

$ws.Range("F79").Value = 14.506060606060601
$ws.Range("G79").Value = 14.3121212121212
$ws.Range("H79").Value = 14.8121212121212

$ws.Range("F81").Value = 59.787500000000001
$ws.Range("G81").Value = 56.612499999999997
$ws.Range("H81").Value = 60.95

$ws.Range("F82").Value = 81.688888888888897
$ws.Range("G82").Value = 84.0833333333334
$ws.Range("H82").Value = 79.5138888888889

$ws.Range("F83").Value = 39.4354430379747
$ws.Range("G83").Value = 39.039240506329101
$ws.Range("H83").Value = 39.755696202531702

$ws.Range("G84").Value = 93.657894736842096

$ws.Range("F87").Value = 74.931578947368394
$ws.Range("G87").Value = 75.515789473684194
$ws.Range("H87").Value = 74.710526315789494

$ws.Range("F88").Value = 50.774999999999999
$ws.Range("G88").Value = 50.125
$ws.Range("H88").Value = 51.1

$ws.Range("F89").Value = 50.010344827586202
$ws.Range("G89").Value = 48.041379310344901
$ws.Range("H89").Value = 50.696551724137898

$ws.Range("F90").Value = 13.8222222222222
$ws.Range("G90").Value = 13.2055555555556
$ws.Range("H90").Value = 14.3722222222222

$ws.Range("G91").Value = 92.196296296296296
$ws.Range("H91").Value = 86.044444444444494

$ws.Range("F94").Value = 58.733333333333398
$ws.Range("G94").Value = 55.608333333333398
$ws.Range("H94").Value = 60.983333333333398

$ws.Range("F95").Value = 85.871428571428595
$ws.Range("G95").Value = 88.55
$ws.Range("H95").Value = 83.821428571428598

$ws.Range("F96").Value = 62.88
$ws.Range("G96").Value = 64.45
$ws.Range("H96").Value = 62.16

$ws.Range("C97").Value = 63.459178571428602
$ws.Range("D97").Value = 57.112178571428601
$ws.Range("E97").Value = 69.979285714285695
$ws.Range("F97").Value = 86.76
$ws.Range("G97").Value = 90.736000000000004
$ws.Range("H97").Value = 83.623999999999995
$ws.Range("I97").Value = 71.535909965103002
$ws.Range("J97").Value = 26.380963914956599
$ws.Range("K97").Value = 2.0831189855276802
$ws.Range("L97").Value = 55.117184691789099
$ws.Range("M97").Value = 16.4187252733138

$ws.Range("C98").Value = 59.288333333333298
$ws.Range("D98").Value = 47.862400000000001
$ws.Range("E98").Value = 70.927466666666703
$ws.Range("F98").Value = 82.144444444444503
$ws.Range("G98").Value = 83.8
$ws.Range("H98").Value = 81.366666666666703
$ws.Range("I98").Value = 50.848997692994601
$ws.Range("J98").Value = 46.324985011071199
$ws.Range("K98").Value = 2.8260202317800598
$ws.Range("L98").Value = 35.533886280230703
$ws.Range("M98").Value = 15.315111412763899

$wb.Save()
